$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 9706
$ws1.Range("F6").Value = 658
$ws1.Range("F7").Value = 103
$ws1.Range("F10").Value = 426
$ws1.Range("F14").Value = 469
$ws1.Range("F15").Value = 12336
$ws1.Range("F18").Value = 313
$ws1.Range("F19").Value = 92
$ws1.Range("F27").Value = 164
$ws1.Range("F28").Value = 2737
$ws1.Range("F29").Value = 49
$ws1.Range("F31").Value = 2107
$ws1.Range("F32").Value = 80
$ws1.Range("F34").Value = 2154
$ws1.Range("F36").Value = 4216
$ws1.Range("F37").Value = 3695
$ws1.Range("F38").Value = 643
$ws1.Range("F39").Value = 2633
$ws1.Range("F44").Value = 570
$ws1.Range("F46").Value = 143
$ws1.Range("F47").Value = 250
$ws1.Range("F49").Value = 136
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 46
$ws2.Range("F13").Value = 43
$ws2.Range("F14").Value = 39
$ws2.Range("F19").Value = 10
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 9706
$ws4.Range("F6").Value = 658
$ws4.Range("F7").Value = 46
$ws4.Range("F8").Value = 103
$ws4.Range("F11").Value = 426
$ws4.Range("F14").Value = 469
$ws4.Range("F15").Value = 12336
$ws4.Range("F17").Value = 313
$ws4.Range("F18").Value = 92
$ws4.Range("F27").Value = 164
$ws4.Range("F28").Value = 2737
$ws4.Range("F29").Value = 2107
$ws4.Range("F30").Value = 80
$ws4.Range("F31").Value = 2154
$ws4.Range("F35").Value = 10
$ws4.Range("F36").Value = 4216
$ws4.Range("F37").Value = 3695
$ws4.Range("F38").Value = 643
$ws4.Range("F39").Value = 2633
$ws4.Range("F44").Value = 570
$ws4.Range("F46").Value = 143
$ws4.Range("F47").Value = 250
$ws4.Range("F49").Value = 136
